$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- 1. "Ativacao:" date changes from 01/01/2012 to 01/01/2023 ---
# Note: cells B15/C15 ("Programa:" row) already happen to reference the very
# same text as B8/C8 ("Ativacao:" row) in the source data, so they must be
# updated too so they keep sharing the (new) value, matching the original
# quirk in the workbook.
foreach ($addr in @("B8", "C8", "B15", "C15")) {
    $ws.Range($addr).Value = "'01/01/2023"
}

$ws.Range("B7").Copy() | Out-Null
$ws.Range("B8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B15").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("C7").Copy() | Out-Null
$ws.Range("C8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C15").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# --- 2. Row 11 ("Objectives:") gains a B/C paragraph ---
$ws.Range("B11").Value = "Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics."
$ws.Range("C11").Value = "Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics."

$ws.Range("B13").Copy() | Out-Null
$ws.Range("B11").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C11").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# --- 3. Row 14 ("Short syllabus:") gains a B/C value ---
$ws.Range("B14").Value = "To be defined, according to the programmed topic."
$ws.Range("C14").Value = "To be defined, according to the programmed topic."

$ws.Range("B13").Copy() | Out-Null
$ws.Range("B14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C14").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# --- 4. Row 16 ("Syllabus:") gains a B/C paragraph ---
$ws.Range("B16").Value = "The content of this elective course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."
$ws.Range("C16").Value = "The content of this elective course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."

$ws.Range("B15").Copy() | Out-Null
$ws.Range("B16").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C16").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false
